$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = 2
$ws.Range("F8").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F18").Value = 3
$ws.Range("F21").Value = 2
$ws.Range("F23").Value = 1
$ws.Range("F25").Value = 1
$ws.Range("F29").Value = 2
$ws.Range("F34").Value = -1
$ws.Range("F37").Value = 2
$ws.Range("F41").Value = 0
$ws.Range("F46").Value = 3
$ws.Range("F48").Value = 3
$ws.Range("F64").Value = 1
$ws.Range("F68").Value = 2
$ws.Range("F69").Value = -3
$ws.Range("F70").Value = 4
$ws.Range("F75").Value = -2
$ws.Range("F76").Value = -5
